$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking" values corrected
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total" values corrected
$ws.Range("B12").Value = 64
$ws.Range("C12").Value = -8
$ws.Range("E12").Value = "56 / 112"
